$wb = $excel.ActiveWorkbook

$wsLider = $wb.Worksheets.Item("Restricciones_del_lider")
$wsLider.Range("A2").Value = "0.8 - x"
$wsLider.Range("B2").Value = "-1.8"
$wsLider.Range("D2").Value = "0.74"
$wsLider.Range("A3").Value = "-0.8 + x"
$wsLider.Range("B3").Value = "-0.19999999999999996"
$wsLider.Range("D3").Value = "0.96"

$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")
$wsFollower.Range("A2").Value = "1.85 - y"
$wsFollower.Range("B2").Value = "-2.85"
$wsFollower.Range("D2").Value = "0.76"
$wsFollower.Range("E2").Value = "-4.1"
$wsFollower.Range("F2").Value = "-9.9"
$wsFollower.Range("A3").Value = "-1.85 + y"
$wsFollower.Range("B3").Value = "0.8500000000000001"
$wsFollower.Range("D3").Value = "0.2"
$wsFollower.Range("E3").Value = "-8.5"
$wsFollower.Range("F3").Value = "-2.6"

$wsPunto = $wb.Worksheets.Item("Punto_modificado")
$wsPunto.Range("A2").Value = "0.8"
$wsPunto.Range("B2").Value = "1.85"

$wsBf = $wb.Worksheets.Item("Vector_bf")
$wsBf.Range("A2").Value = "-9.14325"

$wsBF = $wb.Worksheets.Item("Vector_BF")
$wsBF.Range("A2").Value = "-1.07"
$wsBF.Range("A3").Value = "-33.4"
